$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.1349983333333333
$ws.Cells.Item(2, 8).Value = 0.404995
$ws.Cells.Item(2, 9).Value = 0.06188478316908706
$ws.Cells.Item(2, 10).Value = 0.06188478316908706
$ws.Cells.Item(2, 13).Value = 0.1809866666666667
$ws.Cells.Item(2, 14).Value = 0.54296
$ws.Cells.Item(2, 15).Value = 0.03987407676082905
$ws.Cells.Item(2, 16).Value = 0.03987407676082905
$ws.Cells.Item(2, 17).Value = 0.02443289835555556
$ws.Cells.Item(2, 18).Value = 0.2198960852
$ws.Cells.Item(2, 19).Value = 0.002467598594411439
$ws.Cells.Item(2, 20).Value = 0.002467598594411439

$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.1349983333333333
$ws.Cells.Item(3, 8).Value = 0.404995
$ws.Cells.Item(3, 9).Value = 0.06188478316908706
$ws.Cells.Item(3, 10).Value = 0.06188478316908706
$ws.Cells.Item(3, 15).Value = 0.1057193993302571
$ws.Cells.Item(3, 16).Value = 0.1057193993302571
$ws.Cells.Item(3, 17).Value = 0.06477971524055556
$ws.Cells.Item(3, 18).Value = 0.583017437165
$ws.Cells.Item(3, 19).Value = 0.006542422104319088
$ws.Cells.Item(3, 20).Value = 0.006542422104319088

$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.1349983333333333
$ws.Cells.Item(4, 8).Value = 0.404995
$ws.Cells.Item(4, 9).Value = 0.06188478316908706
$ws.Cells.Item(4, 10).Value = 0.06188478316908706
$ws.Cells.Item(4, 13).Value = 3.878113333333333
$ws.Cells.Item(4, 14).Value = 11.63434
$ws.Cells.Item(4, 15).Value = 0.8544065239089139
$ws.Cells.Item(4, 16).Value = 0.8544065239089139
$ws.Cells.Item(4, 17).Value = 0.5235388364777778
$ws.Cells.Item(4, 18).Value = 4.7118495283
$ws.Cells.Item(4, 19).Value = 0.05287476247035653
$ws.Cells.Item(4, 20).Value = 0.05287476247035653

$ws.Cells.Item(5, 9).Value = 0.4284959871424753
$ws.Cells.Item(5, 10).Value = 0.4284959871424753
$ws.Cells.Item(5, 13).Value = 0.1809866666666667
$ws.Cells.Item(5, 14).Value = 0.54296
$ws.Cells.Item(5, 15).Value = 0.03987407676082905
$ws.Cells.Item(5, 16).Value = 0.03987407676082905
$ws.Cells.Item(5, 17).Value = 0.1691756577866667
$ws.Cells.Item(5, 18).Value = 1.52258092008
$ws.Cells.Item(5, 19).Value = 0.01708588188302628
$ws.Cells.Item(5, 20).Value = 0.01708588188302628

$ws.Cells.Item(6, 9).Value = 0.4284959871424753
$ws.Cells.Item(6, 10).Value = 0.4284959871424753
$ws.Cells.Item(6, 15).Value = 0.1057193993302571
$ws.Cells.Item(6, 16).Value = 0.1057193993302571
$ws.Cells.Item(6, 19).Value = 0.04530033837612806
$ws.Cells.Item(6, 20).Value = 0.04530033837612806

$ws.Cells.Item(7, 9).Value = 0.4284959871424753
$ws.Cells.Item(7, 10).Value = 0.4284959871424753
$ws.Cells.Item(7, 13).Value = 3.878113333333333
$ws.Cells.Item(7, 14).Value = 11.63434
$ws.Cells.Item(7, 15).Value = 0.8544065239089139
$ws.Cells.Item(7, 16).Value = 0.8544065239089139
$ws.Cells.Item(7, 17).Value = 3.625031535313333
$ws.Cells.Item(7, 18).Value = 32.62528381782
$ws.Cells.Item(7, 19).Value = 0.366109766883321
$ws.Cells.Item(7, 20).Value = 0.366109766883321

$ws.Cells.Item(8, 7).Value = 1.111707
$ws.Cells.Item(8, 8).Value = 3.335121
$ws.Cells.Item(8, 9).Value = 0.5096192296884376
$ws.Cells.Item(8, 10).Value = 0.5096192296884376
$ws.Cells.Item(8, 13).Value = 0.1809866666666667
$ws.Cells.Item(8, 14).Value = 0.54296
$ws.Cells.Item(8, 15).Value = 0.03987407676082905
$ws.Cells.Item(8, 16).Value = 0.03987407676082905
$ws.Cells.Item(8, 17).Value = 0.20120414424
$ws.Cells.Item(8, 18).Value = 1.81083729816
$ws.Cells.Item(8, 19).Value = 0.02032059628339133
$ws.Cells.Item(8, 20).Value = 0.02032059628339133

$ws.Cells.Item(9, 7).Value = 1.111707
$ws.Cells.Item(9, 8).Value = 3.335121
$ws.Cells.Item(9, 9).Value = 0.5096192296884376
$ws.Cells.Item(9, 10).Value = 0.5096192296884376
$ws.Cells.Item(9, 15).Value = 0.1057193993302571
$ws.Cells.Item(9, 16).Value = 0.1057193993302571
$ws.Cells.Item(9, 17).Value = 0.533458903623
$ws.Cells.Item(9, 18).Value = 4.801130132607
$ws.Cells.Item(9, 19).Value = 0.05387663884980996
$ws.Cells.Item(9, 20).Value = 0.05387663884980996

$ws.Cells.Item(10, 7).Value = 1.111707
$ws.Cells.Item(10, 8).Value = 3.335121
$ws.Cells.Item(10, 9).Value = 0.5096192296884376
$ws.Cells.Item(10, 10).Value = 0.5096192296884376
$ws.Cells.Item(10, 13).Value = 3.878113333333333
$ws.Cells.Item(10, 14).Value = 11.63434
$ws.Cells.Item(10, 15).Value = 0.8544065239089139
$ws.Cells.Item(10, 16).Value = 0.8544065239089139
$ws.Cells.Item(10, 17).Value = 4.31132573946
$ws.Cells.Item(10, 18).Value = 38.80193165514
$ws.Cells.Item(10, 19).Value = 0.4354219945552364
$ws.Cells.Item(10, 20).Value = 0.4354219945552364
